# edit.ps1 - apply the two content edits described by the diff:
#   1. "Giver the above tables solve the following queries."
#        -> split into three runs "Give" / "n" / " the above tables solve
#           the following queries." (fixing the typo Giver -> Given) while
#           keeping the exact same run formatting (rPr) on every run.
#   2. Wrap the "ANS." run in <w:proofErr w:type="gramStart"/> /
#      <w:proofErr w:type="gramEnd"/> markers.

$d = $word.ActiveDocument

# Shared run formatting (rPr) used throughout this part of the document.
$rpr = '<w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

# -------------------------------------------------------------------------
# Edit 1: "Giver the above tables solve the following queries." -> split
# into three runs, correcting "Giver" to "Given".
# -------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute("Giver the above tables solve the following queries.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    # Re-materialise a plain Range from the Find hit's Start/End before
    # calling InsertXML - calling InsertXML directly on the Range object
    # that Find.Execute mutated leaves the old run behind (duplicated
    # content) instead of replacing it.
    $target1 = $d.Range($rng1.Start, $rng1.End)
    $xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        '<w:r>' + $rpr + '<w:t>Give</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t>n</w:t></w:r>' +
        '<w:r>' + $rpr + '<w:t xml:space="preserve"> the above tables solve the following queries.</w:t></w:r>' +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target1.InsertXML($xml1)
}

# -------------------------------------------------------------------------
# Edit 2: wrap the "ANS." run with gramStart/gramEnd proofErr markers.
# -------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("ANS.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    # Same re-materialisation as above - avoids the stale/duplicated-run
    # InsertXML behaviour on a Find-mutated Range.
    $target2 = $d.Range($rng2.Start, $rng2.End)
    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r>' + $rpr + '<w:t>ANS.</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target2.InsertXML($xml2)
}
